$d = $word.ActiveDocument

$d.Content.Find.Execute("99×90=8910", $true, $false, $false, $false, $false, $true, 1, $false, "14×58=812", 2) | Out-Null
$d.Content.Find.Execute("25×13=325", $true, $false, $false, $false, $false, $true, 1, $false, "44×62=2728", 2) | Out-Null
$d.Content.Find.Execute("13×43=559", $true, $false, $false, $false, $false, $true, 1, $false, "34×18=612", 2) | Out-Null
$d.Content.Find.Execute("80×13=1040", $true, $false, $false, $false, $false, $true, 1, $false, "39×28=1092", 2) | Out-Null
$d.Content.Find.Execute("81×79=6399", $true, $false, $false, $false, $false, $true, 1, $false, "95×41=3895", 2) | Out-Null
$d.Content.Find.Execute("56×17=952", $true, $false, $false, $false, $false, $true, 1, $false, "30×73=2190", 2) | Out-Null
$d.Content.Find.Execute("22×94=2068", $true, $false, $false, $false, $false, $true, 1, $false, "25×66=1650", 2) | Out-Null
$d.Content.Find.Execute("95×25=2375", $true, $false, $false, $false, $false, $true, 1, $false, "46×80=3680", 2) | Out-Null
$d.Content.Find.Execute("93×27=2511", $true, $false, $false, $false, $false, $true, 1, $false, "83×17=1411", 2) | Out-Null
$d.Content.Find.Execute("99×38=3762", $true, $false, $false, $false, $false, $true, 1, $false, "53×21=1113", 2) | Out-Null
$d.Content.Find.Execute("60×53=3180", $true, $false, $false, $false, $false, $true, 1, $false, "43×33=1419", 2) | Out-Null
$d.Content.Find.Execute("30×59=1770", $true, $false, $false, $false, $false, $true, 1, $false, "78×68=5304", 2) | Out-Null
$d.Content.Find.Execute("39×16=624", $true, $false, $false, $false, $false, $true, 1, $false, "65×27=1755", 2) | Out-Null
$d.Content.Find.Execute("64×30=1920", $true, $false, $false, $false, $false, $true, 1, $false, "52×60=3120", 2) | Out-Null
$d.Content.Find.Execute("71×43=3053", $true, $false, $false, $false, $false, $true, 1, $false, "89×54=4806", 2) | Out-Null
$d.Content.Find.Execute("20×16=320", $true, $false, $false, $false, $false, $true, 1, $false, "89×17=1513", 2) | Out-Null
$d.Content.Find.Execute("48×15=720", $true, $false, $false, $false, $false, $true, 1, $false, "76×20=1520", 2) | Out-Null
$d.Content.Find.Execute("18×12=216", $true, $false, $false, $false, $false, $true, 1, $false, "64×21=1344", 2) | Out-Null
$d.Content.Find.Execute("32×41=1312", $true, $false, $false, $false, $false, $true, 1, $false, "44×30=1320", 2) | Out-Null
$d.Content.Find.Execute("16×80=1280", $true, $false, $false, $false, $false, $true, 1, $false, "25×45=1125", 2) | Out-Null
$d.Content.Find.Execute("68×14=952", $true, $false, $false, $false, $false, $true, 1, $false, "48×72=3456", 2) | Out-Null
$d.Content.Find.Execute("89×16=1424", $true, $false, $false, $false, $false, $true, 1, $false, "41×69=2829", 2) | Out-Null
$d.Content.Find.Execute("46×78=3588", $true, $false, $false, $false, $false, $true, 1, $false, "43×56=2408", 2) | Out-Null
$d.Content.Find.Execute("37×54=1998", $true, $false, $false, $false, $false, $true, 1, $false, "44×26=1144", 2) | Out-Null
$d.Content.Find.Execute("70×55=3850", $true, $false, $false, $false, $false, $true, 1, $false, "83×65=5395", 2) | Out-Null
